$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "169.17") need to be
# forced to Text format first, otherwise Excel auto-converts them to a
# floating point number (introducing rounding noise) when .Value is set.
function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

$ws.Range('D2').Value = '66.757.73'
$ws.Range('E2').Value = '  +0.70%  '

$ws.Range('D3').Value = '3.496.97'
$ws.Range('E3').Value = '  +0.25%  '

$ws.Range('E4').Value = '  +0.00%  '

Set-TextValue 'D5' '593.98'
$ws.Range('E5').Value = '  -0.31%  '

Set-TextValue 'D6' '169.17'
$ws.Range('E6').Value = '  -0.29%  '

$ws.Range('E7').Value = '  -0.03%  '

Set-TextValue 'D8' '0.589'
$ws.Range('E8').Value = '  +2.01%  '

Set-TextValue 'D9' '0.133'
$ws.Range('E9').Value = '  +7.00%  '

Set-TextValue 'D10' '7.32'
$ws.Range('E10').Value = '  +0.55%  '

Set-TextValue 'D11' '0.433'
$ws.Range('E11').Value = '  -0.68%  '

$ws.Range('D12').Value = '4.107.72'
$ws.Range('E12').Value = '  +0.36%  '

$ws.Range('E13').Value = '  -0.27%  '

Set-TextValue 'D14' '28.20'
$ws.Range('E14').Value = '  +1.25%  '

Set-TextValue 'D15' '0.0000181'
$ws.Range('E15').Value = '  +1.45%  '

$ws.Range('D16').Value = '66.815.72'
$ws.Range('E16').Value = '  +0.75%  '

$ws.Range('D17').Value = '3.497.80'
$ws.Range('E17').Value = '  -0.01%  '

Set-TextValue 'D18' '6.32'
$ws.Range('E18').Value = '  +0.54%  '

Set-TextValue 'D19' '14.06'
$ws.Range('E19').Value = '  +0.18%  '

Set-TextValue 'D20' '394.12'
$ws.Range('E20').Value = '  +1.67%  '

Set-TextValue 'D21' '7.92'
$ws.Range('E21').Value = '  -0.96%  '

Set-TextValue 'D22' '73.46'
$ws.Range('E22').Value = '  +0.51%  '

Set-TextValue 'D23' '0.998'
$ws.Range('E23').Value = '  -0.11%  '

Set-TextValue 'D24' '0.535'
$ws.Range('E24').Value = '  +1.55%  '

Set-TextValue 'D25' '0.0000123'
$ws.Range('E25').Value = '  +0.32%  '

Set-TextValue 'D26' '10.19'
$ws.Range('E26').Value = '  +0.56%  '

$ws.Range('E27').Value = '  +0.14%  '

Set-TextValue 'D28' '0.999'
$ws.Range('E28').Value = '  +0.18%  '

Set-TextValue 'D29' '6.28'
$ws.Range('E29').Value = '  -1.64%  '

$ws.Range('E30').Value = '  -0.56%  '

Set-TextValue 'D31' '2.06'
$ws.Range('E31').Value = '  -0.42%  '

Set-TextValue 'D32' '23.93'
$ws.Range('E32').Value = '  +2.06%  '

Set-TextValue 'D33' '7.36'
$ws.Range('E33').Value = '  -0.68%  '

$ws.Range('E34').Value = '  +3.57%  '

Set-TextValue 'D35' '163.78'
$ws.Range('E35').Value = '  +2.09%  '

Set-TextValue 'D36' '0.893'
$ws.Range('E36').Value = '  -1.17%  '

Set-TextValue 'D37' '1.91'
$ws.Range('E37').Value = '  -0.79%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D38' '6.83'
$ws.Range('E38').Value = '  +2.54%  '

$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D39' '4.72'
$ws.Range('E39').Value = '  +3.39%  '

Set-TextValue 'D40' '0.0742'
$ws.Range('E40').Value = '  -0.65%  '

Set-TextValue 'D41' '26.32'
$ws.Range('E41').Value = '  -0.15%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.821.10'
$ws.Range('E42').Value = '  +0.52%  '

$ws.Range('B43').Value = 'InjectiveProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D43' '26.93'
$ws.Range('E43').Value = '  -1.26%  '

$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D44' '2.61'
$ws.Range('E44').Value = '  +4.89%  '

Set-TextValue 'D45' '42.75'
$ws.Range('E45').Value = '  -1.38%  '

Set-TextValue 'D46' '0.0311'
$ws.Range('E46').Value = '  -0.46%  '

Set-TextValue 'D47' '341.67'
$ws.Range('E47').Value = '  -2.09%  '

Set-TextValue 'D48' '1.10'
$ws.Range('E48').Value = '  +0.99%  '

Set-TextValue 'D49' '33.70'
$ws.Range('E49').Value = '  +3.27%  '

Set-TextValue 'D50' '0.851'
$ws.Range('E50').Value = '  +0.20%  '

Set-TextValue 'D51' '6.49'
$ws.Range('E51').Value = '  +0.74%  '
